# Updates cryptos list price (D) and 1h volume/change (E) columns
# to the latest scraped values, as inline text cells (matching the
# original workbook's inlineStr cell layout: no numeric coercion,
# no leftover cell styling).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.002.32"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.64%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.913.20"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.56%  "

# Row 4 (percentage only)
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.85%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "598.11"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.41%  "

# Row 7 (percentage only)
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.50%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.198"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.45%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.912.11"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.48%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.425"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +14.39%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.161"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.80%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.90"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.448.34"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.51%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.794.05"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.34%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000192"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.38"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.72%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.899.48"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.03%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.92"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.46%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.77"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.50"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.72%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.32"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.57%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.45%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.33"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.77%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.058.98"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.53%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.20"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.24%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.71"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.82%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000109"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.30%  "

# Row 30 (percentage only)
$ws.Range("E30").Value = "  -0.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.41"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.62%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "502.80"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.16%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.73"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.53%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.24%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.09%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.24"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.22"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.50%  "

# Row 38 (percentage only)
$ws.Range("E38").Value = "  +2.11%  "

# Row 39 (percentage only)
$ws.Range("E39").Value = "  -5.92%  "

# Row 40 (percentage only)
$ws.Range("E40").Value = "  -0.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "179.96"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.88%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.343"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.80%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.01"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.29%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.66"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.95%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0921"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.52%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.16"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.16%  "

# Row 47 (percentage only)
$ws.Range("E47").Value = "  -3.48%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.34"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.17%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.577"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.81%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.663"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.54%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.73"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.76%  "
